# Added solo mode and added up level up functionality
# Fill in the two new time-log entries (rows 19 and 20) on the "Time Log" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Time Log")

# Row 19: 2015-04-29, 16:15 - 17:00, Coding for collisions with enemies / Code
$ws.Range("C19").Value = 42123
$ws.Range("D19").Value = 0.67708333333333337
$ws.Range("E19").Value = 0.70833333333333337
$ws.Range("H19").Value = "Coding for collisions with enemies"
$ws.Range("I19").Value = "Code"

# Row 20: 2015-04-29, 17:00 - 17:15, Testing to make sure code works / Test
$ws.Range("C20").Value = 42123
$ws.Range("D20").Value = 0.70833333333333337
$ws.Range("E20").Value = 0.71875
$ws.Range("H20").Value = "Testing to make sure code works "
$ws.Range("I20").Value = "Test"

# Move the active selection to J19 (matches the author's last recorded selection)
$ws.Range("J19").Select() | Out-Null
